$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.964.32"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.623.08"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D5").Value = "'213.75"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("E9").Value = "  -3.06%  "
$ws.Range("D10").Value = "'18.09"
$ws.Range("E10").Value = "  -7.62%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.847.89"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "1.636.79"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").Value = "25.945.58"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "0.0₃0737"
$ws.Range("E17").Value = "  -3.23%  "
$ws.Range("D18").Value = "'61.14"
$ws.Range("E18").Value = "  -3.52%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "'189.63"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").Value = "'9.54"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "'143.52"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "'1.76"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").Value = "'3.12"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("E33").Value = "  -5.70%  "
$ws.Range("D34").Value = "'2.40"
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").Value = "1.125.03"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "'0.843"
$ws.Range("E37").Value = "  -6.66%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").Value = "'97.55"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "1.759.14"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "'5.17"
$ws.Range("E44").Value = "  -5.36%  "
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").Value = "'54.39"
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("D47").Value = "'0.0523"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").Value = "'7.45"
$ws.Range("E51").Value = "  -3.81%  "
